# Apply crypto price/volume updates per commit
# "Updated cryptos list on Fri Oct 11 17:53:41 UTC 2024 with GitHub Actions"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'62.755.70"
$ws.Range("E2").Value = "  +4.49%  "

$ws.Range("D3").Value = "'2.412.44"
$ws.Range("E3").Value = "  +1.21%  "

$ws.Range("E4").Value = "  +0.31%  "

$ws.Range("D5").Value = "'574.29"
$ws.Range("E5").Value = "  +2.72%  "

$ws.Range("D6").Value = "'144.99"
$ws.Range("E6").Value = "  +5.25%  "

$ws.Range("E7").Value = "  -0.33%  "

$ws.Range("E8").Value = "  +0.64%  "

$ws.Range("D9").Value = "'2.438.66"
$ws.Range("E9").Value = "  +2.25%  "

$ws.Range("E10").Value = "  +5.57%  "

$ws.Range("E11").Value = "  +0.87%  "

$ws.Range("E12").Value = "  +3.83%  "

$ws.Range("E13").Value = "  +4.69%  "

$ws.Range("E14").Value = "  +5.21%  "

$ws.Range("D15").Value = "'0.0000179"
$ws.Range("E15").Value = "  +8.73%  "

$ws.Range("D17").Value = "'62.106.12"
$ws.Range("E17").Value = "  +2.85%  "

$ws.Range("D18").Value = "'2.442.70"
$ws.Range("E18").Value = "  +2.05%  "

$ws.Range("D19").Value = "'7.98"
$ws.Range("E19").Value = "  -5.06%  "

$ws.Range("D20").Value = "'10.91"
$ws.Range("E20").Value = "  +3.41%  "

$ws.Range("D21").Value = "'326.34"
$ws.Range("E21").Value = "  +0.95%  "

$ws.Range("E22").Value = "  +2.94%  "

$ws.Range("E23").Value = "  +15.07%  "

$ws.Range("E24").Value = "  -0.04%  "

$ws.Range("D25").Value = "'65.67"
$ws.Range("E25").Value = "  +2.04%  "

$ws.Range("D26").Value = "'616.30"
$ws.Range("E26").Value = "  +12.32%  "

$ws.Range("D27").Value = "'8.35"
$ws.Range("E27").Value = "  +5.58%  "

$ws.Range("D28").Value = "'0.0₃0988"
$ws.Range("E28").Value = "  +10.08%  "

$ws.Range("D29").Value = "'2.535.04"
$ws.Range("E29").Value = "  +1.23%  "

$ws.Range("D30").Value = "'0.991"
$ws.Range("E30").Value = "  -1.28%  "

$ws.Range("E31").Value = "  +2.81%  "

$ws.Range("E32").Value = "  +10.08%  "

$ws.Range("D33").Value = "'0.138"
$ws.Range("E33").Value = "  +6.20%  "

$ws.Range("E34").Value = "  +2.97%  "

$ws.Range("E35").Value = "  +5.63%  "

$ws.Range("E36").Value = "  -0.47%  "

$ws.Range("E37").Value = "  +6.14%  "

$ws.Range("B38").Value = "Monero"
$ws.Range("C38").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D38").Value = "'152.85"
$ws.Range("E38").Value = "  +0.05%  "

$ws.Range("B39").Value = "PolygonEcosystemToken"
$ws.Range("C39").Value = "https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol"
$ws.Range("D39").Value = "'0.373"
$ws.Range("E39").Value = "  +2.01%  "

$ws.Range("B40").Value = "RenderToken"
$ws.Range("C40").Value = "https://coinranking.com/coin/vfo5XYwcV+rendertoken-render"
$ws.Range("D40").Value = "'5.38"
$ws.Range("E40").Value = "  +8.38%  "

$ws.Range("D41").Value = "'18.58"
$ws.Range("E41").Value = "  +1.95%  "

$ws.Range("D42").Value = "'2.77"
$ws.Range("E42").Value = "  +21.73%  "

$ws.Range("E43").Value = "  +8.47%  "

$ws.Range("D44").Value = "'42.28"
$ws.Range("E44").Value = "  +2.87%  "

$ws.Range("E46").Value = "  +0.80%  "

$ws.Range("D47").Value = "'144.51"
$ws.Range("E47").Value = "  +2.20%  "

$ws.Range("D48").Value = "'3.58"
$ws.Range("E48").Value = "  +2.81%  "

$ws.Range("D49").Value = "'20.20"
$ws.Range("E49").Value = "  +7.74%  "

$ws.Range("D50").Value = "'0.600"
$ws.Range("E50").Value = "  +2.43%  "

$ws.Range("D51").Value = "'0.0514"
$ws.Range("E51").Value = "  +3.67%  "

